# Atualização de bases das ligas, do dia: 19-04-2024 às 21:40
#
# 1) Rows 83 and 84 had their match data (everything except the "id" in
#    column A) swapped.
# 2) A new row (93) was appended with the next match fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the contents of rows 83 and 84 (columns B..AC). The "id"
#    values in column A (81 / 82) stay where they are. Values are set
#    as literals (rather than read from one row and copied to the
#    other) so that full floating point precision is preserved.
# ---------------------------------------------------------------------
$ws.Cells.Item(83, 2).Value = 7301364
$ws.Cells.Item(83, 3).Value = "Canada Premier League"
$ws.Cells.Item(83, 4).Value = "Canada Premier League"
$ws.Cells.Item(83, 5).Value = 45206.75
$ws.Cells.Item(83, 6).Value = "Forge FC"
$ws.Cells.Item(83, 7).Value = "Atletico Ottawa"
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 1
$ws.Cells.Item(83, 10).Value = "A"
$ws.Cells.Item(83, 11).Value = 1.8
$ws.Cells.Item(83, 12).Value = 3.6
$ws.Cells.Item(83, 13).Value = 3.5
$ws.Cells.Item(83, 14).Value = 1.533
$ws.Cells.Item(83, 15).Value = 3.8
$ws.Cells.Item(83, 16).Value = 5
$ws.Cells.Item(83, 17).Value = -1
$ws.Cells.Item(83, 18).Value = 1.975
$ws.Cells.Item(83, 19).Value = 1.825
$ws.Cells.Item(83, 20).Value = 2.5
$ws.Cells.Item(83, 21).Value = 1.9
$ws.Cells.Item(83, 22).Value = 1.9
$ws.Cells.Item(83, 23).Value = -1
$ws.Cells.Item(83, 24).Value = -1
$ws.Cells.Item(83, 25).Value = 4
$ws.Cells.Item(83, 26).Value = -1
$ws.Cells.Item(83, 27).Value = 0.825
$ws.Cells.Item(83, 28).Value = -1
$ws.Cells.Item(83, 29).Value = 0.8999999999999999

$ws.Cells.Item(84, 2).Value = 6227884
$ws.Cells.Item(84, 3).Value = "Canada Premier League"
$ws.Cells.Item(84, 4).Value = "Canada Premier League"
$ws.Cells.Item(84, 5).Value = 45206.75
$ws.Cells.Item(84, 6).Value = "Cavalry FC"
$ws.Cells.Item(84, 7).Value = "Pacific FC CA"
$ws.Cells.Item(84, 8).Value = 3
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = "H"
$ws.Cells.Item(84, 11).Value = 2.25
$ws.Cells.Item(84, 12).Value = 3.1
$ws.Cells.Item(84, 13).Value = 2.875
$ws.Cells.Item(84, 14).Value = 2.05
$ws.Cells.Item(84, 15).Value = 3.2
$ws.Cells.Item(84, 16).Value = 3.2
$ws.Cells.Item(84, 17).Value = -0.25
$ws.Cells.Item(84, 18).Value = 1.825
$ws.Cells.Item(84, 19).Value = 1.975
$ws.Cells.Item(84, 20).Value = 2.5
$ws.Cells.Item(84, 21).Value = 1.825
$ws.Cells.Item(84, 22).Value = 1.975
$ws.Cells.Item(84, 23).Value = 1.05
$ws.Cells.Item(84, 24).Value = -1
$ws.Cells.Item(84, 25).Value = -1
$ws.Cells.Item(84, 26).Value = 0.825
$ws.Cells.Item(84, 27).Value = -1
$ws.Cells.Item(84, 28).Value = 0.825
$ws.Cells.Item(84, 29).Value = -1

# ---------------------------------------------------------------------
# 2) Append new row 93 with the fresh fixture data. Copy the number
#    formatting used for the "id" column (A) and the "Date" column (E)
#    from the previous row so the new row matches the existing style.
# ---------------------------------------------------------------------
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E92").Copy()
$ws.Range("E93").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(93, 1).Value = 91
$ws.Cells.Item(93, 2).Value = 7802935
$ws.Cells.Item(93, 3).Value = "Canada Premier League"
$ws.Cells.Item(93, 4).Value = "Canada Premier League"
$ws.Cells.Item(93, 5).Value = 45401.95833333334
$ws.Cells.Item(93, 6).Value = "Pacific FC CA"
$ws.Cells.Item(93, 7).Value = "Valour FC"
$ws.Cells.Item(93, 11).Value = 1.727
$ws.Cells.Item(93, 12).Value = 3.5
$ws.Cells.Item(93, 13).Value = 4
$ws.Cells.Item(93, 14).Value = 1.533
$ws.Cells.Item(93, 15).Value = 4.2
$ws.Cells.Item(93, 16).Value = 4.5
$ws.Cells.Item(93, 17).Value = -1
$ws.Cells.Item(93, 18).Value = 1.95
$ws.Cells.Item(93, 19).Value = 1.85
$ws.Cells.Item(93, 20).Value = 2.5
$ws.Cells.Item(93, 21).Value = 1.925
$ws.Cells.Item(93, 22).Value = 1.875
$ws.Cells.Item(93, 23).Value = 0
$ws.Cells.Item(93, 24).Value = 0
$ws.Cells.Item(93, 25).Value = 0
$ws.Cells.Item(93, 26).Value = 0
$ws.Cells.Item(93, 27).Value = 0

# Columns H, I, J (FTHG, FTAG, FTR) and AB, AC (PL_AhOver, PL_AhUnder)
# are intentionally left blank for this not-yet-played fixture.

Write-Output "Edit applied."
